$wb = $excel.ActiveWorkbook

# ---------- Sheet 1: LP1912 ----------
$ws1 = $wb.Worksheets.Item("LP1912")

# Header updates
$ws1.Cells.Item(2,1).Value = 'Última actualización: 10:50:37'
$ws1.Cells.Item(3,1).Value = 'Total filas: 186'

# Swap rows 61/62 (A, C, D columns; B and E stay the same)
$ws1.Cells.Item(61,1).Value = '05:47:32'
$ws1.Cells.Item(61,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(61,4).Value = 105
$ws1.Cells.Item(62,1).Value = '06:02:16'
$ws1.Cells.Item(62,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(62,4).Value = 90

# Swap Linea (C) for rows 109/110
$ws1.Cells.Item(109,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(110,3).Value = '17_ROMERO'

# Swap Linea (C) for rows 120/121
$ws1.Cells.Item(120,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(121,3).Value = '23_HERNANDEZ'

# Rebuild rows 152-191 (re-sorted/expanded schedule tail)
$ws1.Cells.Item(152,1).Value = '10:50:37'
$ws1.Cells.Item(152,2).Value = '10:51'
$ws1.Cells.Item(152,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(152,4).Value = 1
$ws1.Cells.Item(152,5).Value = 'LP1912'
$ws1.Cells.Item(153,1).Value = '10:11:11'
$ws1.Cells.Item(153,2).Value = '10:52'
$ws1.Cells.Item(153,3).Value = '15_ABASTO'
$ws1.Cells.Item(153,4).Value = 41
$ws1.Cells.Item(153,5).Value = 'LP1912'
$ws1.Cells.Item(154,1).Value = '09:25:30'
$ws1.Cells.Item(154,2).Value = '10:53'
$ws1.Cells.Item(154,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(154,4).Value = 88
$ws1.Cells.Item(154,5).Value = 'LP1912'
$ws1.Cells.Item(155,1).Value = '10:50:37'
$ws1.Cells.Item(155,2).Value = '10:56'
$ws1.Cells.Item(155,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(155,4).Value = 6
$ws1.Cells.Item(155,5).Value = 'LP1912'
$ws1.Cells.Item(156,1).Value = '10:11:11'
$ws1.Cells.Item(156,2).Value = '10:57'
$ws1.Cells.Item(156,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(156,4).Value = 46
$ws1.Cells.Item(156,5).Value = 'LP1912'
$ws1.Cells.Item(157,1).Value = '10:50:37'
$ws1.Cells.Item(157,2).Value = '11:01'
$ws1.Cells.Item(157,3).Value = '215C_EL PATO'
$ws1.Cells.Item(157,4).Value = 11
$ws1.Cells.Item(157,5).Value = 'LP1912'
$ws1.Cells.Item(158,1).Value = '09:25:30'
$ws1.Cells.Item(158,2).Value = '11:02'
$ws1.Cells.Item(158,3).Value = '215C_EL PATO'
$ws1.Cells.Item(158,4).Value = 97
$ws1.Cells.Item(158,5).Value = 'LP1912'
$ws1.Cells.Item(159,1).Value = '10:50:37'
$ws1.Cells.Item(159,2).Value = '11:03'
$ws1.Cells.Item(159,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(159,4).Value = 13
$ws1.Cells.Item(159,5).Value = 'LP1912'
$ws1.Cells.Item(160,1).Value = '10:11:11'
$ws1.Cells.Item(160,2).Value = '11:04'
$ws1.Cells.Item(160,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(160,4).Value = 53
$ws1.Cells.Item(160,5).Value = 'LP1912'
$ws1.Cells.Item(161,1).Value = '10:50:37'
$ws1.Cells.Item(161,2).Value = '11:04'
$ws1.Cells.Item(161,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(161,4).Value = 14
$ws1.Cells.Item(161,5).Value = 'LP1912'
$ws1.Cells.Item(162,1).Value = '10:11:11'
$ws1.Cells.Item(162,2).Value = '11:05'
$ws1.Cells.Item(162,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(162,4).Value = 54
$ws1.Cells.Item(162,5).Value = 'LP1912'
$ws1.Cells.Item(163,1).Value = '09:25:30'
$ws1.Cells.Item(163,2).Value = '11:06'
$ws1.Cells.Item(163,3).Value = '16_P MOR-167 Y 521'
$ws1.Cells.Item(163,4).Value = 101
$ws1.Cells.Item(163,5).Value = 'LP1912'
$ws1.Cells.Item(164,1).Value = '10:11:11'
$ws1.Cells.Item(164,2).Value = '11:07'
$ws1.Cells.Item(164,3).Value = '16_P MOR-167 Y 521'
$ws1.Cells.Item(164,4).Value = 56
$ws1.Cells.Item(164,5).Value = 'LP1912'
$ws1.Cells.Item(165,1).Value = '10:11:11'
$ws1.Cells.Item(165,2).Value = '11:11'
$ws1.Cells.Item(165,3).Value = '10_OLMOS'
$ws1.Cells.Item(165,4).Value = 60
$ws1.Cells.Item(165,5).Value = 'LP1912'
$ws1.Cells.Item(166,1).Value = '10:11:11'
$ws1.Cells.Item(166,2).Value = '11:12'
$ws1.Cells.Item(166,3).Value = '15_ABASTO'
$ws1.Cells.Item(166,4).Value = 61
$ws1.Cells.Item(166,5).Value = 'LP1912'
$ws1.Cells.Item(167,1).Value = '09:25:30'
$ws1.Cells.Item(167,2).Value = '11:19'
$ws1.Cells.Item(167,3).Value = '86_EST CHICA-ESC AGRARIA'
$ws1.Cells.Item(167,4).Value = 114
$ws1.Cells.Item(167,5).Value = 'LP1912'
$ws1.Cells.Item(168,1).Value = '10:11:11'
$ws1.Cells.Item(168,2).Value = '11:20'
$ws1.Cells.Item(168,3).Value = '86_EST CHICA-ESC AGRARIA'
$ws1.Cells.Item(168,4).Value = 69
$ws1.Cells.Item(168,5).Value = 'LP1912'
$ws1.Cells.Item(169,1).Value = '09:25:30'
$ws1.Cells.Item(169,2).Value = '11:21'
$ws1.Cells.Item(169,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(169,4).Value = 116
$ws1.Cells.Item(169,5).Value = 'LP1912'
$ws1.Cells.Item(170,1).Value = '10:11:11'
$ws1.Cells.Item(170,2).Value = '11:27'
$ws1.Cells.Item(170,3).Value = '225_C ROCA-H SUR'
$ws1.Cells.Item(170,4).Value = 76
$ws1.Cells.Item(170,5).Value = 'LP1912'
$ws1.Cells.Item(171,1).Value = '10:11:11'
$ws1.Cells.Item(171,2).Value = '11:32'
$ws1.Cells.Item(171,3).Value = '81_EL PELIGRO'
$ws1.Cells.Item(171,4).Value = 81
$ws1.Cells.Item(171,5).Value = 'LP1912'
$ws1.Cells.Item(172,1).Value = '10:50:37'
$ws1.Cells.Item(172,2).Value = '11:34'
$ws1.Cells.Item(172,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(172,4).Value = 44
$ws1.Cells.Item(172,5).Value = 'LP1912'
$ws1.Cells.Item(173,1).Value = '10:50:37'
$ws1.Cells.Item(173,2).Value = '11:35'
$ws1.Cells.Item(173,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(173,4).Value = 45
$ws1.Cells.Item(173,5).Value = 'LP1912'
$ws1.Cells.Item(174,1).Value = '10:11:11'
$ws1.Cells.Item(174,2).Value = '11:38'
$ws1.Cells.Item(174,3).Value = '10_OLMOS'
$ws1.Cells.Item(174,4).Value = 87
$ws1.Cells.Item(174,5).Value = 'LP1912'
$ws1.Cells.Item(175,1).Value = '10:50:37'
$ws1.Cells.Item(175,2).Value = '11:41'
$ws1.Cells.Item(175,3).Value = '17_ROMERO'
$ws1.Cells.Item(175,4).Value = 51
$ws1.Cells.Item(175,5).Value = 'LP1912'
$ws1.Cells.Item(176,1).Value = '10:11:11'
$ws1.Cells.Item(176,2).Value = '11:42'
$ws1.Cells.Item(176,3).Value = '17_ROMERO'
$ws1.Cells.Item(176,4).Value = 91
$ws1.Cells.Item(176,5).Value = 'LP1912'
$ws1.Cells.Item(177,1).Value = '10:50:37'
$ws1.Cells.Item(177,2).Value = '11:43'
$ws1.Cells.Item(177,3).Value = '10_OLMOS'
$ws1.Cells.Item(177,4).Value = 53
$ws1.Cells.Item(177,5).Value = 'LP1912'
$ws1.Cells.Item(178,1).Value = '10:11:11'
$ws1.Cells.Item(178,2).Value = '11:51'
$ws1.Cells.Item(178,3).Value = '215B_EL PATO'
$ws1.Cells.Item(178,4).Value = 100
$ws1.Cells.Item(178,5).Value = 'LP1912'
$ws1.Cells.Item(179,1).Value = '10:11:11'
$ws1.Cells.Item(179,2).Value = '11:59'
$ws1.Cells.Item(179,3).Value = '225_GOMEZ'
$ws1.Cells.Item(179,4).Value = 108
$ws1.Cells.Item(179,5).Value = 'LP1912'
$ws1.Cells.Item(180,1).Value = '10:11:11'
$ws1.Cells.Item(180,2).Value = '12:02'
$ws1.Cells.Item(180,3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws1.Cells.Item(180,4).Value = 111
$ws1.Cells.Item(180,5).Value = 'LP1912'
$ws1.Cells.Item(181,1).Value = '10:50:37'
$ws1.Cells.Item(181,2).Value = '12:06'
$ws1.Cells.Item(181,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(181,4).Value = 76
$ws1.Cells.Item(181,5).Value = 'LP1912'
$ws1.Cells.Item(182,1).Value = '10:50:37'
$ws1.Cells.Item(182,2).Value = '12:06'
$ws1.Cells.Item(182,3).Value = '14_ABASTO'
$ws1.Cells.Item(182,4).Value = 76
$ws1.Cells.Item(182,5).Value = 'LP1912'
$ws1.Cells.Item(183,1).Value = '10:11:11'
$ws1.Cells.Item(183,2).Value = '12:07'
$ws1.Cells.Item(183,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(183,4).Value = 116
$ws1.Cells.Item(183,5).Value = 'LP1912'
$ws1.Cells.Item(184,1).Value = '10:50:37'
$ws1.Cells.Item(184,2).Value = '12:07'
$ws1.Cells.Item(184,3).Value = '10_OLMOS'
$ws1.Cells.Item(184,4).Value = 77
$ws1.Cells.Item(184,5).Value = 'LP1912'
$ws1.Cells.Item(185,1).Value = '10:11:11'
$ws1.Cells.Item(185,2).Value = '12:07'
$ws1.Cells.Item(185,3).Value = '14_ABASTO'
$ws1.Cells.Item(185,4).Value = 116
$ws1.Cells.Item(185,5).Value = 'LP1912'
$ws1.Cells.Item(186,1).Value = '10:50:37'
$ws1.Cells.Item(186,2).Value = '12:20'
$ws1.Cells.Item(186,3).Value = '215A_EL PATO'
$ws1.Cells.Item(186,4).Value = 90
$ws1.Cells.Item(186,5).Value = 'LP1912'
$ws1.Cells.Item(187,1).Value = '10:50:37'
$ws1.Cells.Item(187,2).Value = '12:21'
$ws1.Cells.Item(187,3).Value = '14_ABASTO'
$ws1.Cells.Item(187,4).Value = 91
$ws1.Cells.Item(187,5).Value = 'LP1912'
$ws1.Cells.Item(188,1).Value = '10:50:37'
$ws1.Cells.Item(188,2).Value = '12:21'
$ws1.Cells.Item(188,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(188,4).Value = 91
$ws1.Cells.Item(188,5).Value = 'LP1912'
$ws1.Cells.Item(189,1).Value = '10:50:37'
$ws1.Cells.Item(189,2).Value = '12:22'
$ws1.Cells.Item(189,3).Value = '17_ROMERO'
$ws1.Cells.Item(189,4).Value = 92
$ws1.Cells.Item(189,5).Value = 'LP1912'
$ws1.Cells.Item(190,1).Value = '10:50:37'
$ws1.Cells.Item(190,2).Value = '12:36'
$ws1.Cells.Item(190,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(190,4).Value = 106
$ws1.Cells.Item(190,5).Value = 'LP1912'
$ws1.Cells.Item(191,1).Value = '10:50:37'
$ws1.Cells.Item(191,2).Value = '12:38'
$ws1.Cells.Item(191,3).Value = '17_179 Y 38'
$ws1.Cells.Item(191,4).Value = 108
$ws1.Cells.Item(191,5).Value = 'LP1912'

# ---------- Sheet 2: LP1912-215 ----------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2,1).Value = 'Última actualización: 10:50:37'
$ws2.Cells.Item(3,1).Value = 'Total filas: 24'

# Rebuild rows 26-29
$ws2.Cells.Item(26,1).Value = '10:50:37'
$ws2.Cells.Item(26,2).Value = '11:01'
$ws2.Cells.Item(26,3).Value = '215C_EL PATO'
$ws2.Cells.Item(26,4).Value = 11
$ws2.Cells.Item(26,5).Value = 'LP1912'
$ws2.Cells.Item(27,1).Value = '09:25:30'
$ws2.Cells.Item(27,2).Value = '11:02'
$ws2.Cells.Item(27,3).Value = '215C_EL PATO'
$ws2.Cells.Item(27,4).Value = 97
$ws2.Cells.Item(27,5).Value = 'LP1912'
$ws2.Cells.Item(28,1).Value = '10:11:11'
$ws2.Cells.Item(28,2).Value = '11:51'
$ws2.Cells.Item(28,3).Value = '215B_EL PATO'
$ws2.Cells.Item(28,4).Value = 100
$ws2.Cells.Item(28,5).Value = 'LP1912'
$ws2.Cells.Item(29,1).Value = '10:50:37'
$ws2.Cells.Item(29,2).Value = '12:20'
$ws2.Cells.Item(29,3).Value = '215A_EL PATO'
$ws2.Cells.Item(29,4).Value = 90
$ws2.Cells.Item(29,5).Value = 'LP1912'

# ---------- Sheet 3: 6203-6173 ----------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2,1).Value = 'Última actualización: 10:50:37'
$ws3.Cells.Item(3,1).Value = 'Total filas: 30'

# Rebuild rows 31-35
$ws3.Cells.Item(31,1).Value = '10:50:37'
$ws3.Cells.Item(31,2).Value = '10:53'
$ws3.Cells.Item(31,3).Value = '215A_LA PLATA'
$ws3.Cells.Item(31,4).Value = 3
$ws3.Cells.Item(31,5).Value = 'L6173'
$ws3.Cells.Item(32,1).Value = '09:25:30'
$ws3.Cells.Item(32,2).Value = '10:54'
$ws3.Cells.Item(32,3).Value = '215A_LA PLATA'
$ws3.Cells.Item(32,4).Value = 89
$ws3.Cells.Item(32,5).Value = 'L6173'
$ws3.Cells.Item(33,1).Value = '10:50:37'
$ws3.Cells.Item(33,2).Value = '11:13'
$ws3.Cells.Item(33,3).Value = '215C_LA PLATA'
$ws3.Cells.Item(33,4).Value = 23
$ws3.Cells.Item(33,5).Value = 'L6203'
$ws3.Cells.Item(34,1).Value = '09:25:30'
$ws3.Cells.Item(34,2).Value = '11:14'
$ws3.Cells.Item(34,3).Value = '215C_LA PLATA'
$ws3.Cells.Item(34,4).Value = 109
$ws3.Cells.Item(34,5).Value = 'L6203'
$ws3.Cells.Item(35,1).Value = '10:11:11'
$ws3.Cells.Item(35,2).Value = '12:04'
$ws3.Cells.Item(35,3).Value = '215A_LA PLATA'
$ws3.Cells.Item(35,4).Value = 113
$ws3.Cells.Item(35,5).Value = 'L6173'
